$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Fix the category label on row 19 (was pointing at the stale
#    "LIVE, BILLBOARD" text, should read "LIVE, BILLBOARD, TRAFFIC"
#    just like row 18).
# ------------------------------------------------------------------
$ws.Range("A19").Value = "LIVE, BILLBOARD, TRAFFIC"

# ------------------------------------------------------------------
# 2) Add the two new Calgary, Canada rows (20 and 21) underneath the
#    existing table, reusing the formatting already used by the
#    table body (row 19) so no stray new styles are introduced.
# ------------------------------------------------------------------

# --- Row 20: "Calgary Downtown View" ---
$ws.Range("B19").Copy() | Out-Null
$ws.Range("B20:C20").PasteSpecial(-4122) | Out-Null
$ws.Range("I19").Copy() | Out-Null
$ws.Range("I20").PasteSpecial(-4122) | Out-Null
$ws.Range("J19").Copy() | Out-Null
$ws.Range("J20").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D20:F20").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null

$ws.Range("A20").Value = "LIVE, DOWNTOWN, RIVER, TRAFFIC, BUILDING"
$ws.Range("D20").Value = "Calgary Downtown View"
$ws.Range("E20").Value = "Calgary"
$ws.Range("F20").Value = "Canada"
$ws.Range("G20").Value = "MwcqP3ta6RI"
$ws.Range("I20").Value = 51.051458907264099
$ws.Range("J20").Value = -114.05785799826501
$ws.Range("B20").Formula = "=I20+(K20*0.000001)-(L20*0.000001)"
$ws.Range("C20").Formula = "=J20+(M20*0.0001)-(N20*0.0001)"

# --- Row 21: "Central Memorial Park" ---
$ws.Range("B19").Copy() | Out-Null
$ws.Range("B21:C21").PasteSpecial(-4122) | Out-Null
$ws.Range("I19").Copy() | Out-Null
$ws.Range("I21").PasteSpecial(-4122) | Out-Null
$ws.Range("J19").Copy() | Out-Null
$ws.Range("J21").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").Copy() | Out-Null
$ws.Range("A21:F21").PasteSpecial(-4122) | Out-Null

$ws.Range("A21").Value = "LIVE, PARK"
$ws.Range("D21").Value = "Central Memorial Park"
$ws.Range("E21").Value = "Calgary"
$ws.Range("F21").Value = "Canada"
$ws.Range("G21").Value = "xsRDTfuksyI"
$ws.Range("I21").Value = 51.0411974223206
$ws.Range("J21").Value = -114.070227502294
$ws.Range("B21").Formula = "=I21+(K21*0.000001)-(L21*0.000001)"
$ws.Range("C21").Formula = "=J21+(M21*0.0001)-(N21*0.0001)"

# ------------------------------------------------------------------
# 3) Update the selection to reflect where the editor ended up.
# ------------------------------------------------------------------
$ws.Range("A20").Select() | Out-Null
